# Apply crypto price/volume updates per commit diff (Thu Jun 13 18:48:20 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells must stay plain text (values like "0.0890", "79.06",
# "1.00" etc. would otherwise be re-interpreted as numbers and lose their exact
# textual representation), so force text format before assigning them.
$dCells = @("D2","D3","D5","D6","D7","D11","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D26","D28","D29","D34","D35","D36","D38","D40","D41","D42","D43","D45","D46","D47","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.816.22'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '3.471.12'
$ws.Range("E3").Value = '  -3.29%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '603.97'
$ws.Range("E5").Value = '  -3.43%  '
$ws.Range("D6").Value = '148.26'
$ws.Range("E6").Value = '  -6.09%  '
$ws.Range("D7").Value = '3.468.93'
$ws.Range("E7").Value = '  -3.17%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").Value = '7.58'
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("E12").Value = '  -3.61%  '
$ws.Range("E13").Value = '  -4.66%  '
$ws.Range("D14").Value = '31.84'
$ws.Range("E14").Value = '  -5.32%  '
$ws.Range("D15").Value = '4.053.33'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").Value = '3.468.32'
$ws.Range("E16").Value = '  -3.13%  '
$ws.Range("D17").Value = '66.829.77'
$ws.Range("E17").Value = '  -4.07%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  -4.89%  '
$ws.Range("D20").Value = '15.43'
$ws.Range("E20").Value = '  -4.55%  '
$ws.Range("D21").Value = '10.12'
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D22").Value = '440.81'
$ws.Range("E22").Value = '  -4.68%  '
$ws.Range("D23").Value = '0.613'
$ws.Range("E23").Value = '  -5.13%  '
$ws.Range("D24").Value = '79.06'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '3.608.13'
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("E27").Value = '  -9.98%  '
$ws.Range("D28").Value = '9.80'
$ws.Range("E28").Value = '  -8.34%  '
$ws.Range("D29").Value = '8.40'
$ws.Range("E29").Value = '  -9.21%  '
$ws.Range("E30").Value = '  -5.55%  '
$ws.Range("E31").Value = '  -7.31%  '
$ws.Range("E32").Value = '  -3.09%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").Value = '25.45'
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").Value = '6.07'
$ws.Range("E35").Value = '  -6.89%  '
$ws.Range("D36").Value = '3.460.44'
$ws.Range("E36").Value = '  -3.56%  '
$ws.Range("E37").Value = '  -6.90%  '
$ws.Range("D38").Value = '7.94'
$ws.Range("E38").Value = '  -5.51%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '175.06'
$ws.Range("E41").Value = '  -2.68%  '
$ws.Range("D42").Value = '0.0890'
$ws.Range("E42").Value = '  -3.78%  '
$ws.Range("D43").Value = '2.15'
$ws.Range("E43").Value = '  -10.74%  '
$ws.Range("E44").Value = '  -4.41%  '
$ws.Range("D45").Value = '0.887'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").Value = '29.28'
$ws.Range("E46").Value = '  -6.40%  '
$ws.Range("D47").Value = '46.19'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  -9.93%  '
$ws.Range("E49").Value = '  -10.26%  '
$ws.Range("D50").Value = '7.47'
$ws.Range("E50").Value = '  -4.60%  '
$ws.Range("D51").Value = '0.989'
$ws.Range("E51").Value = '  -4.82%  '
